$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'27.162.79"
$c.ClearFormats()
$ws.Range("E2").Value = "  -0.99%  "

$c = $ws.Range("D3")
$c.Value = "'1.781.39"
$c.ClearFormats()
$ws.Range("E3").Value = "  +1.00%  "

$c = $ws.Range("D4")
$c.Value = "'1.011"
$c.ClearFormats()
$ws.Range("E4").Value = "  +1.59%  "

$c = $ws.Range("D5")
$c.Value = "'335.90"
$c.ClearFormats()
$ws.Range("E5").Value = "  +0.15%  "

$c = $ws.Range("D6")
$c.Value = "'1.005"
$c.ClearFormats()
$ws.Range("E6").Value = "  +1.38%  "

$c = $ws.Range("D7")
$c.Value = "'0.3785"
$c.ClearFormats()
$ws.Range("E7").Value = "  +0.45%  "

$c = $ws.Range("D8")
$c.Value = "'0.3429"
$c.ClearFormats()
$ws.Range("E8").Value = "  -0.23%  "

$c = $ws.Range("D9")
$c.Value = "'48.54"
$c.ClearFormats()
$ws.Range("E9").Value = "  -0.89%  "

$c = $ws.Range("D10")
$c.Value = "'1.191"
$c.ClearFormats()
$ws.Range("E10").Value = "  -1.34%  "

$c = $ws.Range("D11")
$c.Value = "'0.07451"
$c.ClearFormats()
$ws.Range("E11").Value = "  -2.15%  "

$c = $ws.Range("D12")
$c.Value = "'1.005"
$c.ClearFormats()
$ws.Range("E12").Value = "  +1.11%  "

$c = $ws.Range("D13")
$c.Value = "'21.81"
$c.ClearFormats()
$ws.Range("E13").Value = "  +3.48%  "

$c = $ws.Range("D14")
$c.Value = "'6.431"
$c.ClearFormats()
$ws.Range("E14").Value = "  -0.80%  "

$c = $ws.Range("D15")
$c.Value = "'1.789.55"
$c.ClearFormats()
$ws.Range("E15").Value = "  +1.72%  "

$c = $ws.Range("D16")
$c.Value = "'7.054"
$c.ClearFormats()
$ws.Range("E16").Value = "  -1.14%  "

$c = $ws.Range("D17")
$c.Value = "'0.00001095"
$c.ClearFormats()
$ws.Range("E17").Value = "  -0.60%  "

$c = $ws.Range("D18")
$c.Value = "'0.06681"
$c.ClearFormats()
$ws.Range("E18").Value = "  -0.81%  "

$c = $ws.Range("D19")
$c.Value = "'84.43"
$c.ClearFormats()
$ws.Range("E19").Value = "  +0.33%  "

$c = $ws.Range("D20")
$c.Value = "'1.003"
$c.ClearFormats()
$ws.Range("E20").Value = "  +1.07%  "

$c = $ws.Range("D21")
$c.Value = "'6.523"
$c.ClearFormats()
$ws.Range("E21").Value = "  +3.60%  "

$c = $ws.Range("D22")
$c.Value = "'17.29"
$c.ClearFormats()
$ws.Range("E22").Value = "  +0.39%  "

$c = $ws.Range("D23")
$c.Value = "'27.230.97"
$c.ClearFormats()
$ws.Range("E23").Value = "  -0.61%  "

$c = $ws.Range("D24")
$c.Value = "'12.44"
$c.ClearFormats()
$ws.Range("E24").Value = "  -4.39%  "

$c = $ws.Range("D25")
$c.Value = "'2.433"
$c.ClearFormats()
$ws.Range("E25").Value = "  -0.87%  "

$c = $ws.Range("D26")
$c.Value = "'1.497"
$c.ClearFormats()
$ws.Range("E26").Value = "  -0.51%  "

$c = $ws.Range("D27")
$c.Value = "'2.542"
$c.ClearFormats()
$ws.Range("E27").Value = "  +3.04%  "

$c = $ws.Range("D28")
$c.Value = "'21.42"
$c.ClearFormats()
$ws.Range("E28").Value = "  +7.46%  "

$c = $ws.Range("D29")
$c.Value = "'152.83"
$c.ClearFormats()
$ws.Range("E29").Value = "  -0.15%  "

$c = $ws.Range("D30")
$c.Value = "'1.993.23"
$c.ClearFormats()
$ws.Range("E30").Value = "  +1.93%  "

$c = $ws.Range("D31")
$c.Value = "'133.41"
$c.ClearFormats()
$ws.Range("E31").Value = "  -0.67%  "

$c = $ws.Range("D32")
$c.Value = "'4.055"
$c.ClearFormats()
$ws.Range("E32").Value = "  -0.91%  "

$c = $ws.Range("D33")
$c.Value = "'6.027"
$c.ClearFormats()
$ws.Range("E33").Value = "  -2.18%  "

$c = $ws.Range("D34")
$c.Value = "'0.08677"
$c.ClearFormats()
$ws.Range("E34").Value = "  +0.34%  "

$c = $ws.Range("D35")
$c.Value = "'13.07"
$c.ClearFormats()
$ws.Range("E35").Value = "  -0.22%  "

$c = $ws.Range("D36")
$c.Value = "'1.650"
$c.ClearFormats()
$ws.Range("E36").Value = "  -2.55%  "

$c = $ws.Range("D37")
$c.Value = "'5.439"
$c.ClearFormats()
$ws.Range("E37").Value = "  -1.41%  "

$c = $ws.Range("D38")
$c.Value = "'0.6843"
$c.ClearFormats()
$ws.Range("E38").Value = "  +3.99%  "

$c = $ws.Range("D39")
$c.Value = "'0.06373"
$c.ClearFormats()
$ws.Range("E39").Value = "  -0.28%  "

$c = $ws.Range("D40")
$c.Value = "'8.803"
$c.ClearFormats()
$ws.Range("E40").Value = "  +1.97%  "

$c = $ws.Range("D41")
$c.Value = "'0.2186"
$c.ClearFormats()
$ws.Range("E41").Value = "  -1.26%  "

$c = $ws.Range("D42")
$c.Value = "'0.02330"
$c.ClearFormats()
$ws.Range("E42").Value = "  -2.13%  "

$c = $ws.Range("D43")
$c.Value = "'1.260"
$c.ClearFormats()
$ws.Range("E43").Value = "  +2.53%  "

$c = $ws.Range("D44")
$c.Value = "'14.54"
$c.ClearFormats()
$ws.Range("E44").Value = "  +0.45%  "

$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$c = $ws.Range("D45")
$c.Value = "'1.003"
$c.ClearFormats()
$ws.Range("E45").Value = "  +1.10%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$c = $ws.Range("D46")
$c.Value = "'0.6408"
$c.ClearFormats()
$ws.Range("E46").Value = "  +1.05%  "

$c = $ws.Range("D47")
$c.Value = "'3.851"
$c.ClearFormats()
$ws.Range("E47").Value = "  -2.44%  "

$c = $ws.Range("D48")
$c.Value = "'2.124"
$c.ClearFormats()
$ws.Range("E48").Value = "  +0.55%  "

$c = $ws.Range("D49")
$c.Value = "'129.05"
$c.ClearFormats()
$ws.Range("E49").Value = "  -1.04%  "

$c = $ws.Range("D50")
$c.Value = "'0.07185"
$c.ClearFormats()
$ws.Range("E50").Value = "  -1.73%  "

$c = $ws.Range("D51")
$c.Value = "'79.28"
$c.ClearFormats()
$ws.Range("E51").Value = "  +0.11%  "
